$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab (workbook.xml <sheet name="...">) from "Through 2022-07-15"
# to "Through 2022-07-16".
$ws.Name = "Through 2022-07-16"

# Update the header label for the current-month column (B1), reflecting the new
# "through" date.
$ws.Range("B1").Value = "July 2022 (through July 16)"

# Updated/added carjacking counts for various neighborhoods/months (new data
# for 2022-07-16, plus a handful of historical corrections in the same
# "July of year X" columns).
$ws.Range("B2").Value = 11    # Austin, July 2022
$ws.Range("P2").Value = 4     # Austin, July 2020
$ws.Range("AK2").Value = 4    # Austin, July 2017

$ws.Range("AK4").Value = 1    # Auburn Gresham, July 2017

$ws.Range("I5").Value = 2     # Garfield Park, July 2021
$ws.Range("P5").Value = 7     # Garfield Park, July 2020
$ws.Range("AD5").Value = 3    # Garfield Park, July 2018

$ws.Range("P8").Value = 11    # North Lawndale, July 2020

$ws.Range("I12").Value = 2    # Hyde Park, July 2021

$ws.Range("P15").Value = 2    # Douglas, July 2020

$ws.Range("P19").Value = 2    # South Shore, July 2020
$ws.Range("W19").Value = 3    # South Shore, July 2019

$ws.Range("B24").Value = 1    # South Deering, July 2022

$ws.Range("B26").Value = 4    # Little Village, July 2022

$ws.Range("B27").Value = 4    # Lincoln Park, July 2022

$ws.Range("B38").Value = 1    # West Town, July 2022

$ws.Range("P39").Value = 1    # Wicker Park, July 2020

$ws.Range("P52").Value = 8    # Chatham, July 2020
$ws.Range("AD52").Value = 2   # Chatham, July 2018

$ws.Range("AK57").Value = 1   # Woodlawn, July 2017

$ws.Range("B70").Value = 3    # Edgewater, July 2022

$ws.Range("B78").Value = 2    # Lake View, July 2022

$ws.Range("B85").Value = 1    # North Center, July 2022

$ws.Range("AK92").Value = 1   # Sauganash,Forest Glen, July 2017
